$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1143.5883
$ws.Range("I28").Value = 1176.0667
$ws.Range("K28").Value = 1176.0667
$ws.Range("M28").Value = -691.0667000000001
$ws.Range("H32").Value = 2189.9092
$ws.Range("J32").Value = 2328.9
$ws.Range("L32").Value = 2328.9
$ws.Range("N32").Value = -2980.9
$ws.Range("H75").Value = 31500
$ws.Range("J75").Value = 31500
$ws.Range("L75").Value = 31500
$ws.Range("N75").Value = -33372
$ws.Range("H78").Value = 31500
$ws.Range("J78").Value = 31500
$ws.Range("L78").Value = 94500
$ws.Range("N78").Value = -103860
$ws.Range("H112").Value = 2473.1667
$ws.Range("J112").Value = 2373
$ws.Range("L112").Value = 7119
$ws.Range("N112").Value = -9335
$ws.Range("H132").Value = 806.8570999999999
$ws.Range("I132").Value = 874.6667
$ws.Range("K132").Value = 2624.0001
$ws.Range("M132").Value = -94.0001000000002
$ws.Range("H135").Value = 22.333334
$ws.Range("I135").Value = 22.333334
$ws.Range("K135").Value = 201.000006
$ws.Range("M135").Value = 2333.999994
$ws.Range("H137").Value = 2181.5454
$ws.Range("I137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("M137").Value = -3450
$ws.Range("H138").Value = 2181.5642
$ws.Range("J138").Value = 2826.4
$ws.Range("L138").Value = 8479.200000000001
$ws.Range("N138").Value = -18759.2
$ws.Range("H141").Value = 1605.091
$ws.Range("I141").Value = 1586.2858
$ws.Range("K141").Value = 4758.857400000001
$ws.Range("M141").Value = 421.1425999999992

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3743.5588
$ws.Range("I32").Value = 3743.5588
$ws.Range("K32").Value = 3743.5588
$ws.Range("M32").Value = -3456.5588
$ws.Range("H74").Value = 2664.818
$ws.Range("I74").Value = 2681.3
$ws.Range("K74").Value = 2681.3
$ws.Range("M74").Value = -1807.3
$ws.Range("H77").Value = 2664.818
$ws.Range("I77").Value = 2681.3
$ws.Range("K77").Value = 13406.5
$ws.Range("M77").Value = -9038.5
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802
$ws.Range("H97").Value = 439.6154
$ws.Range("I97").Value = 518.6842
$ws.Range("J97").Value = 225
$ws.Range("K97").Value = 518.6842
$ws.Range("L97").Value = 225
$ws.Range("M97").Value = -22.68420000000003
$ws.Range("N97").Value = -1217

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1480.6
$ws.Range("I94").Value = 1480.6
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1480.6
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1029.6
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 1692.2222
$ws.Range("I105").Value = 1711.4
$ws.Range("K105").Value = 1711.4
$ws.Range("M105").Value = 35.59999999999991
$ws.Range("H107").Value = 3834.6428
$ws.Range("I107").Value = 3473.75
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 3473.75
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -1553.75
$ws.Range("N107").Value = -9840
$ws.Range("H134").Value = 1299.6666
$ws.Range("I134").Value = 1362.75
$ws.Range("K134").Value = 4088.25
$ws.Range("M134").Value = -1553.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 633.5833
$ws.Range("I5").Value = 262
$ws.Range("J5").Value = 1376.75
$ws.Range("K5").Value = 262
$ws.Range("L5").Value = 1376.75
$ws.Range("M5").Value = -150
$ws.Range("N5").Value = -1600.75
$ws.Range("H16").Value = 5865.1113
$ws.Range("I16").Value = 5700.3335
$ws.Range("J16").Value = 6194.6665
$ws.Range("K16").Value = 5700.3335
$ws.Range("L16").Value = 6194.6665
$ws.Range("M16").Value = -5413.3335
$ws.Range("N16").Value = -6768.6665
$ws.Range("H58").Value = 1358.409
$ws.Range("I58").Value = 1363.9474
$ws.Range("K58").Value = 1363.9474
$ws.Range("M58").Value = -1160.9474
$ws.Range("H105").Value = 2470.6875
$ws.Range("I105").Value = 1976.1666
$ws.Range("K105").Value = 1976.1666
$ws.Range("M105").Value = -229.1666
$ws.Range("H113").Value = 5865.1113
$ws.Range("I113").Value = 5700.3335
$ws.Range("J113").Value = 6194.6665
$ws.Range("K113").Value = 5700.3335
$ws.Range("L113").Value = 6194.6665
$ws.Range("M113").Value = -3530.3335
$ws.Range("N113").Value = -10534.6665
$ws.Range("H132").Value = 2351.8462
$ws.Range("I132").Value = 2434
$ws.Range("K132").Value = 7302
$ws.Range("M132").Value = -4772
$ws.Range("H134").Value = 1249.25
$ws.Range("I134").Value = 1249.25
$ws.Range("K134").Value = 3747.75
$ws.Range("M134").Value = -1212.75
$ws.Range("H136").Value = 1358.409
$ws.Range("I136").Value = 1363.9474
$ws.Range("K136").Value = 4091.8422
$ws.Range("M136").Value = -1541.8422
$ws.Range("H141").Value = 161997.33
$ws.Range("I141").Value = 43991
$ws.Range("J141").Value = 176748.12
$ws.Range("K141").Value = 43991
$ws.Range("L141").Value = 176748.12
$ws.Range("M141").Value = -38811
$ws.Range("N141").Value = -187108.12

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 18650
$ws.Range("I117").Value = 633.3333
$ws.Range("J117").Value = 36666.668
$ws.Range("K117").Value = 1899.9999
$ws.Range("L117").Value = 110000.004
$ws.Range("M117").Value = 1542.0001
$ws.Range("N117").Value = -116884.004
$ws.Range("H129").Value = 1113593.1
$ws.Range("I129").Value = 561
$ws.Range("J129").Value = 2504883.2
$ws.Range("K129").Value = 1683
$ws.Range("L129").Value = 7514649.600000001
$ws.Range("M129").Value = 3317
$ws.Range("N129").Value = -7524649.600000001
$ws.Range("H131").Value = 502705.44
$ws.Range("J131").Value = 669694.8
$ws.Range("L131").Value = 2009084.4
$ws.Range("N131").Value = -2019164.4
$ws.Range("H132").Value = 7997
$ws.Range("I132").Value = 7997
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 71973
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -69443
$ws.Range("N132").ClearContents()
$ws.Range("H140").Value = 12000.5
$ws.Range("I140").Value = 3002
$ws.Range("J140").Value = 17726.818
$ws.Range("K140").Value = 9006
$ws.Range("L140").Value = 53180.454
$ws.Range("M140").Value = -3826
$ws.Range("N140").Value = -63540.454
$ws.Range("H141").Value = 8687.25
$ws.Range("I141").Value = 8687.25
$ws.Range("K141").Value = 26061.75
$ws.Range("M141").Value = -20881.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 44999
$ws.Range("J26").Value = 44999
$ws.Range("L26").Value = 44999
$ws.Range("N26").Value = -45559
$ws.Range("H50").Value = 44999
$ws.Range("J50").Value = 44999
$ws.Range("L50").Value = 44999
$ws.Range("N50").Value = -45995
$ws.Range("H53").Value = 240000
$ws.Range("I53").Value = 240000
$ws.Range("K53").Value = 240000
$ws.Range("M53").Value = -239369

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1826.8572
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 2032
$ws.Range("I22").Value = 1957.6923
$ws.Range("K22").Value = 1957.6923
$ws.Range("M22").Value = -1662.6923
$ws.Range("H27").Value = 2032
$ws.Range("I27").Value = 1957.6923
$ws.Range("K27").Value = 1957.6923
$ws.Range("M27").Value = -1850.6923
$ws.Range("H40").Value = 1945.1666
$ws.Range("I40").Value = 1729.25
$ws.Range("J40").Value = 3024.75
$ws.Range("K40").Value = 1729.25
$ws.Range("L40").Value = 3024.75
$ws.Range("M40").Value = -1593.25
$ws.Range("N40").Value = -3296.75
$ws.Range("H132").Value = 2557.7144
$ws.Range("I132").Value = 2149.8333
$ws.Range("K132").Value = 6449.499899999999
$ws.Range("M132").Value = -3919.499899999999
$ws.Range("H136").Value = 3239.3809
$ws.Range("I136").Value = 2779.111
$ws.Range("K136").Value = 8337.332999999999
$ws.Range("M136").Value = -5787.332999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1266.3334
$ws.Range("I122").Value = 1266.3334
$ws.Range("K122").Value = 3799.0002
$ws.Range("M122").Value = -1349.0002
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws.Range("H132").Value = 2672.9048
$ws.Range("I132").Value = 2848.8235
$ws.Range("K132").Value = 8546.470499999999
$ws.Range("M132").Value = -6016.470499999999
